$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Insert a new worksheet ("Sheet1") between "Pinout" and
# "UML StateMachine" - this is where the PID-tuning helper table
# used for "send pid to backend" lives.
# ---------------------------------------------------------------
$pinout = $wb.Worksheets.Item("Pinout")
$new = $wb.Worksheets.Add($null, $pinout)

# ---------------------------------------------------------------
# Header row: three repeated Pk / kI / kD groups (cols A-C, E-G, I-K)
# ---------------------------------------------------------------
$new.Range("A1").Value = "Pk"
$new.Range("B1").Value = "kI"
$new.Range("C1").Value = "kD"

$new.Range("E1").Value = "Pk"
$new.Range("F1").Value = "kI"
$new.Range("G1").Value = "kD"

$new.Range("I1").Value = "Pk"
$new.Range("J1").Value = "kI"
$new.Range("K1").Value = "kD"

# ---------------------------------------------------------------
# Sample gain values under each group (rows 2-7)
# ---------------------------------------------------------------
$vals = @(1, 10, 20, 40, 80, 200)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $row = 2 + $i
    $new.Cells.Item($row, 1).Value = $vals[$i]
    $new.Cells.Item($row, 5).Value = $vals[$i]
    $new.Cells.Item($row, 9).Value = $vals[$i]
}

# ---------------------------------------------------------------
# PID formula notes in column O (written O15 first, then O14, so
# the shared-string table fills in the same order as the source).
# ---------------------------------------------------------------
$new.Range("O15").Value = "Kd = Kp * Td"
$new.Range("O14").Value = "Ki =Kp/Ti"

$new.Rows.Item(14).RowHeight = 25.5
$new.Rows.Item(15).RowHeight = 25.5

foreach ($addr in @("O14", "O15")) {
    $rng = $new.Range($addr)
    $font = $rng.Font
    $font.Bold = $true
    $font.Size = 10
    $font.Name = "Arial Unicode MS"
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
    $rng.IndentLevel = 1
}

# Widen column O so the PID notes are readable.
$new.Columns.Item(15).ColumnWidth = 23

# ---------------------------------------------------------------
# Make the new sheet the active one, with P14 selected, matching
# where the author was working when the file was saved.
# ---------------------------------------------------------------
$new.Activate()
$new.Range("P14").Select()

# ---------------------------------------------------------------
# Pinout sheet: selection moved on from the old C32:C33 to C34.
# ---------------------------------------------------------------
$pinout.Range("C34").Select()
$new.Activate()
$new.Range("P14").Select()
